$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.826.51'
$ws.Range('E2').Value = '  -4.37%  '
$ws.Range('D3').Value = '''3.008.23'
$ws.Range('E3').Value = '  -4.93%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''526.23'
$ws.Range('E5').Value = '  -6.67%  '
$ws.Range('D6').Value = '''127.32'
$ws.Range('E6').Value = '  -10.94%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '''2.998.08'
$ws.Range('E8').Value = '  -5.06%  '
$ws.Range('D9').Value = '''0.489'
$ws.Range('E9').Value = '  -1.57%  '
$ws.Range('E10').Value = '  -4.33%  '
$ws.Range('D11').Value = '''5.94'
$ws.Range('E11').Value = '  -12.55%  '
$ws.Range('D12').Value = '''0.442'
$ws.Range('E13').Value = '  -1.88%  '
$ws.Range('D14').Value = '''32.93'
$ws.Range('E14').Value = '  -10.22%  '
$ws.Range('D15').Value = '''3.489.59'
$ws.Range('D16').Value = '''61.840.30'
$ws.Range('E16').Value = '  -4.48%  '
$ws.Range('E17').Value = '  -2.83%  '
$ws.Range('D18').Value = '''3.018.28'
$ws.Range('E18').Value = '  -4.55%  '
$ws.Range('E19').Value = '  -6.65%  '
$ws.Range('D20').Value = '''467.62'
$ws.Range('D21').Value = '''12.92'
$ws.Range('E21').Value = '  -8.22%  '
$ws.Range('D22').Value = '''0.677'
$ws.Range('E22').Value = '  -5.97%  '
$ws.Range('D23').Value = '''6.84'
$ws.Range('E23').Value = '  -8.47%  '
$ws.Range('D24').Value = '''77.43'
$ws.Range('E24').Value = '  -2.24%  '
$ws.Range('D25').Value = '''11.60'
$ws.Range('E25').Value = '  -9.37%  '
$ws.Range('D26').Value = '''0.997'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('D27').Value = '''2.61'
$ws.Range('E27').Value = '  -8.39%  '
$ws.Range('D28').Value = '''7.84'
$ws.Range('E28').Value = '  -11.19%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').Value = '''25.17'
$ws.Range('E30').Value = '  -5.54%  '
$ws.Range('E31').Value = '  -16.33%  '
$ws.Range('E32').Value = '  -5.73%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = '''55.78'
$ws.Range('E33').Value = '  +2.65%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = '''2.30'
$ws.Range('E34').Value = '  -11.71%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '''5.74'
$ws.Range('E35').Value = '  -6.03%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '''5.04'
$ws.Range('E36').Value = '  -6.38%  '
$ws.Range('D37').Value = '''459.74'
$ws.Range('E37').Value = '  -16.72%  '
$ws.Range('D38').Value = '''3.013.34'
$ws.Range('E38').Value = '  -4.95%  '
$ws.Range('D39').Value = '''0.0382'
$ws.Range('E39').Value = '  -11.89%  '
$ws.Range('D40').Value = '''0.0768'
$ws.Range('E40').Value = '  -6.92%  '
$ws.Range('D41').Value = '''0.110'
$ws.Range('E41').Value = '  -10.56%  '
$ws.Range('D42').Value = '''7.82'
$ws.Range('E42').Value = '  -5.98%  '
$ws.Range('D43').Value = '''2.45'
$ws.Range('E43').Value = '  -12.21%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '''0.241'
$ws.Range('E45').Value = '  -9.47%  '
$ws.Range('E46').Value = '  -12.07%  '
$ws.Range('D47').Value = '''0.0₃0509'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''116.16'
$ws.Range('E48').Value = '  -4.64%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.105'
$ws.Range('E49').Value = '  -3.00%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '''23.38'
$ws.Range('E50').Value = '  -7.26%  '
$ws.Range('B51').Value = 'CoreDAO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D51').Value = '''2.25'
$ws.Range('E51').Value = '  -6.76%  '
